$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(80, 8).Value = 46297052
$ws.Cells.Item(80, 9).Value = 14706373
$ws.Cells.Item(80, 10).Value = 100001200
$ws.Cells.Item(80, 11).Value = 44119119
$ws.Cells.Item(80, 12).Value = 300003600
$ws.Cells.Item(80, 13).Value = -44118121
$ws.Cells.Item(80, 14).Value = -300005596

$ws.Cells.Item(83, 8).Value = 46297052
$ws.Cells.Item(83, 9).Value = 14706373
$ws.Cells.Item(83, 10).Value = 100001200
$ws.Cells.Item(83, 11).Value = 132357357
$ws.Cells.Item(83, 12).Value = 900010800
$ws.Cells.Item(83, 13).Value = -132352365
$ws.Cells.Item(83, 14).Value = -900020784

$ws.Cells.Item(106, 8).Value = 3380.8125
$ws.Cells.Item(106, 9).Value = 3257.3076
$ws.Cells.Item(106, 11).Value = 3257.3076
$ws.Cells.Item(106, 13).Value = -2626.3076

$ws.Cells.Item(125, 8).Value = 2933.9375
$ws.Cells.Item(125, 9).Value = 2117.1
$ws.Cells.Item(125, 10).Value = 4295.3335
$ws.Cells.Item(125, 11).Value = 19053.9
$ws.Cells.Item(125, 12).Value = 38658.0015
$ws.Cells.Item(125, 13).Value = -16593.9
$ws.Cells.Item(125, 14).Value = -43578.0015

$ws.Cells.Item(132, 8).Value = 377370.75
$ws.Cells.Item(132, 9).Value = 471710.4
$ws.Cells.Item(132, 10).Value = 18880
$ws.Cells.Item(132, 11).Value = 1415131.2
$ws.Cells.Item(132, 12).Value = 56640
$ws.Cells.Item(132, 13).Value = -1412601.2
$ws.Cells.Item(132, 14).Value = -61700

$ws.Cells.Item(138, 8).Value = 1905.36
$ws.Cells.Item(138, 9).Value = 1490.6428
$ws.Cells.Item(138, 11).Value = 4471.928400000001
$ws.Cells.Item(138, 13).Value = 668.0715999999993

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 13027.11
$ws.Cells.Item(32, 9).Value = 11889.387
$ws.Cells.Item(32, 11).Value = 11889.387
$ws.Cells.Item(32, 13).Value = -11602.387

$ws.Cells.Item(45, 8).Value = 3114.4285
$ws.Cells.Item(45, 9).Value = 2776.5557
$ws.Cells.Item(45, 10).Value = 3722.6
$ws.Cells.Item(45, 11).Value = 2776.5557
$ws.Cells.Item(45, 12).Value = 3722.6
$ws.Cells.Item(45, 13).Value = -2399.5557
$ws.Cells.Item(45, 14).Value = -4476.6

$ws.Cells.Item(63, 8).Value = 551.6667
$ws.Cells.Item(63, 9).Value = 551.6667
$ws.Cells.Item(63, 11).Value = 551.6667
$ws.Cells.Item(63, 13).Value = 134.3333

$ws.Cells.Item(66, 8).Value = 551.6667
$ws.Cells.Item(66, 9).Value = 551.6667
$ws.Cells.Item(66, 11).Value = 2758.3335
$ws.Cells.Item(66, 13).Value = 673.6665000000003

$ws.Cells.Item(74, 8).Value = 12501676
$ws.Cells.Item(74, 9).Value = 16668027
$ws.Cells.Item(74, 10).Value = 2619.8
$ws.Cells.Item(74, 11).Value = 16668027
$ws.Cells.Item(74, 12).Value = 2619.8
$ws.Cells.Item(74, 13).Value = -16667153
$ws.Cells.Item(74, 14).Value = -4367.8

$ws.Cells.Item(77, 8).Value = 12501676
$ws.Cells.Item(77, 9).Value = 16668027
$ws.Cells.Item(77, 10).Value = 2619.8
$ws.Cells.Item(77, 11).Value = 83340135
$ws.Cells.Item(77, 12).Value = 13099
$ws.Cells.Item(77, 13).Value = -83335767
$ws.Cells.Item(77, 14).Value = -21835

$ws.Cells.Item(122, 8).Value = 2470.2285
$ws.Cells.Item(122, 9).Value = 2018.8
$ws.Cells.Item(122, 10).Value = 5178.8
$ws.Cells.Item(122, 11).Value = 6056.4
$ws.Cells.Item(122, 12).Value = 15536.4
$ws.Cells.Item(122, 13).Value = -3606.4
$ws.Cells.Item(122, 14).Value = -20436.4

$ws.Cells.Item(132, 8).Value = 11984.695
$ws.Cells.Item(132, 9).Value = 14762.25
$ws.Cells.Item(132, 10).Value = 7096.2
$ws.Cells.Item(132, 11).Value = 44286.75
$ws.Cells.Item(132, 12).Value = 21288.6
$ws.Cells.Item(132, 13).Value = -41756.75
$ws.Cells.Item(132, 14).Value = -26348.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 2804.9556
$ws.Cells.Item(86, 9).Value = 2470.8333
$ws.Cells.Item(86, 11).Value = 2470.8333
$ws.Cells.Item(86, 13).Value = -1347.8333

$ws.Cells.Item(89, 8).Value = 2804.9556
$ws.Cells.Item(89, 9).Value = 2470.8333
$ws.Cells.Item(89, 11).Value = 12354.1665
$ws.Cells.Item(89, 13).Value = -6738.166499999999

$ws.Cells.Item(94, 8).Value = 1051.0605
$ws.Cells.Item(94, 9).Value = 1097.12
$ws.Cells.Item(94, 11).Value = 1097.12
$ws.Cells.Item(94, 13).Value = -646.1199999999999

$ws.Cells.Item(96, 8).Value = 56927.375
$ws.Cells.Item(96, 9).Value = 55185
$ws.Cells.Item(96, 10).Value = 57972.8
$ws.Cells.Item(96, 11).Value = 55185
$ws.Cells.Item(96, 12).Value = 57972.8
$ws.Cells.Item(96, 13).Value = -52439
$ws.Cells.Item(96, 14).Value = -63464.8

$ws.Cells.Item(99, 8).Value = 3308.6875
$ws.Cells.Item(99, 9).Value = 3882.8
$ws.Cells.Item(99, 10).Value = 2351.8333
$ws.Cells.Item(99, 11).Value = 3882.8
$ws.Cells.Item(99, 12).Value = 2351.8333
$ws.Cells.Item(99, 13).Value = -2384.8
$ws.Cells.Item(99, 14).Value = -5347.8333

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 39651.883
$ws.Cells.Item(62, 9).Value = 28010
$ws.Cells.Item(62, 10).Value = 52749
$ws.Cells.Item(62, 11).Value = 28010
$ws.Cells.Item(62, 12).Value = 52749
$ws.Cells.Item(62, 13).Value = -27386
$ws.Cells.Item(62, 14).Value = -53997

$ws.Cells.Item(65, 8).Value = 39651.883
$ws.Cells.Item(65, 9).Value = 28010
$ws.Cells.Item(65, 10).Value = 52749
$ws.Cells.Item(65, 11).Value = 140050
$ws.Cells.Item(65, 12).Value = 263745
$ws.Cells.Item(65, 13).Value = -136930
$ws.Cells.Item(65, 14).Value = -269985

$ws.Cells.Item(80, 8).Value = 32564
$ws.Cells.Item(80, 10).Value = 30128
$ws.Cells.Item(80, 12).Value = 30128
$ws.Cells.Item(80, 14).Value = -32374

$ws.Cells.Item(83, 8).Value = 32564
$ws.Cells.Item(83, 10).Value = 30128
$ws.Cells.Item(83, 12).Value = 90384
$ws.Cells.Item(83, 14).Value = -101616

$ws.Cells.Item(134, 8).Value = 2342.4482
$ws.Cells.Item(134, 9).Value = 2293.963
$ws.Cells.Item(134, 10).Value = 2997
$ws.Cells.Item(134, 11).Value = 6881.889000000001
$ws.Cells.Item(134, 12).Value = 8991
$ws.Cells.Item(134, 13).Value = -4346.889000000001
$ws.Cells.Item(134, 14).Value = -14061

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value = 1354.375
$ws.Cells.Item(34, 10).Value = 2865.8333
$ws.Cells.Item(34, 12).Value = 8597.499899999999
$ws.Cells.Item(34, 14).Value = -8765.499899999999

$ws.Cells.Item(55, 8).Value = 5800
$ws.Cells.Item(55, 10).Value = 5800
$ws.Cells.Item(55, 12).Value = 17400
$ws.Cells.Item(55, 14).Value = -17754

$ws.Cells.Item(81, 8).Value = 4141.5
$ws.Cells.Item(81, 9).Value = 2826.4
$ws.Cells.Item(81, 10).Value = 6333.3335
$ws.Cells.Item(81, 11).Value = 8479.200000000001
$ws.Cells.Item(81, 12).Value = 19000.0005
$ws.Cells.Item(81, 13).Value = -7356.200000000001
$ws.Cells.Item(81, 14).Value = -21246.0005

$ws.Cells.Item(84, 8).Value = 4141.5
$ws.Cells.Item(84, 9).Value = 2826.4
$ws.Cells.Item(84, 10).Value = 6333.3335
$ws.Cells.Item(84, 11).Value = 25437.6
$ws.Cells.Item(84, 12).Value = 57000.0015
$ws.Cells.Item(84, 13).Value = -19821.6
$ws.Cells.Item(84, 14).Value = -68232.0015

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(48, 8).Value = 20000
$ws.Cells.Item(48, 9).Value = 20000
$ws.Cells.Item(48, 10).Value = 0
$ws.Cells.Item(48, 11).Value = 20000
$ws.Cells.Item(48, 13).Value = -19515
$ws.Cells.Item(48, 14).ClearContents()

$ws.Cells.Item(126, 8).Value = 1685.6538
$ws.Cells.Item(126, 9).Value = 1558.2646
$ws.Cells.Item(126, 10).Value = 1926.2778
$ws.Cells.Item(126, 11).Value = 4674.793799999999
$ws.Cells.Item(126, 12).Value = 5778.8334
$ws.Cells.Item(126, 13).Value = -2204.793799999999
$ws.Cells.Item(126, 14).Value = -10718.8334

$ws.Cells.Item(132, 8).Value = 88717.25999999999
$ws.Cells.Item(132, 9).Value = 118338.586
$ws.Cells.Item(132, 11).Value = 355015.758
$ws.Cells.Item(132, 13).Value = -352485.758

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 5930.773
$ws.Cells.Item(7, 9).Value = 5093
$ws.Cells.Item(7, 10).Value = 8164.8335
$ws.Cells.Item(7, 11).Value = 5093
$ws.Cells.Item(7, 12).Value = 8164.8335
$ws.Cells.Item(7, 13).Value = -4981
$ws.Cells.Item(7, 14).Value = -8388.833500000001

$ws.Cells.Item(16, 8).Value = 2030.475
$ws.Cells.Item(16, 9).Value = 1585.7812
$ws.Cells.Item(16, 10).Value = 3809.25
$ws.Cells.Item(16, 11).Value = 1585.7812
$ws.Cells.Item(16, 12).Value = 3809.25
$ws.Cells.Item(16, 13).Value = -1415.7812
$ws.Cells.Item(16, 14).Value = -4149.25

$ws.Cells.Item(40, 8).Value = 16671597
$ws.Cells.Item(40, 9).Value = 15154842
$ws.Cells.Item(40, 11).Value = 15154842
$ws.Cells.Item(40, 13).Value = -15154706

$ws.Cells.Item(122, 8).Value = 7979.081
$ws.Cells.Item(122, 9).Value = 4128.1177
$ws.Cells.Item(122, 11).Value = 12384.3531
$ws.Cells.Item(122, 13).Value = -9934.3531

$ws.Cells.Item(126, 8).Value = 5930.773
$ws.Cells.Item(126, 9).Value = 5093
$ws.Cells.Item(126, 10).Value = 8164.8335
$ws.Cells.Item(126, 11).Value = 15279
$ws.Cells.Item(126, 12).Value = 24494.5005
$ws.Cells.Item(126, 13).Value = -12809
$ws.Cells.Item(126, 14).Value = -29434.5005

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(47, 8).Value = 0
$ws.Cells.Item(47, 10).Value = 0
$ws.Cells.Item(47, 14).ClearContents()

$ws.Cells.Item(126, 8).Value = 142858450
